$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Phase 1 - capture every piece of text content that needs to move,
# BEFORE any mutation happens, so offsets/paragraph numbers are still
# the original ("before") ones and nothing has been overwritten yet.
# Paragraph indices below were established by inspecting the document:
#   6  Objetivos            -> plain paragraph  (TEXT_A)
#   7  Objetivos            -> empty italic run  (TEXT_B, paragraph removed later)
#   9  Docente(s)           -> ListBullet, run1 (TEXT_C1) <br/> run2 (TEXT_C2)
#   11 Programa resumido    -> plain paragraph, 2 runs w/ <br/> (TEXT_D)
#   12 Programa resumido    -> italic paragraph (TEXT_E, paragraph removed later)
#   14 Programa             -> plain paragraph (TEXT_F)
#   15 Programa             -> italic paragraph (TEXT_G, stays put)
#   17 Avaliacao            -> ListBullet, Metodo/Criterio/Norma labels + answers
#   19 Bibliografia         -> plain paragraph, 3 runs w/ <br/> (TEXT_K)
# ---------------------------------------------------------------------

$CR = [char]13
$VT = [char]11

$textA = $d.Paragraphs(6).Range.Text.TrimEnd($CR)

$p9Parts = $d.Paragraphs(9).Range.Text.TrimEnd($CR).Split($VT)
$textC1 = $p9Parts[0]
$textC2 = $p9Parts[1]

$textD = $d.Paragraphs(11).Range.Text.TrimEnd($CR)

$textF = $d.Paragraphs(14).Range.Text.TrimEnd($CR)

$p17Text = $d.Paragraphs(17).Range.Text.TrimEnd($CR)
$p17Parts = $p17Text.Split($VT)

$labelMetodo = "M" + [char]0xE9 + "todo: "
$labelCriterio = "Crit" + [char]0xE9 + "rio: "
$labelNorma = "Norma de recupera" + [char]0xE7 + [char]0xE3 + "o: "

$textH = $p17Parts[0].Substring($labelMetodo.Length)
$textI = $p17Parts[1].Substring($labelCriterio.Length)
$textJ = $p17Parts[2].Substring($labelNorma.Length)

$p19Text = $d.Paragraphs(19).Range.Text.TrimEnd($CR)
$textK = $p19Text

# ---------------------------------------------------------------------
# Phase 2 - write the new text into place using Find.Execute scoped to
# each destination paragraph's Range, matching the *old* content of
# that paragraph/run and swapping in the captured value. This only
# touches the matched span, leaving sibling runs (and their own
# formatting) untouched.
# ---------------------------------------------------------------------

# Objetivos paragraph (6): old TEXT_A -> TEXT_D
$r = $d.Paragraphs(6).Range
$r.Find.Execute($textA, $true, $false, $false, $false, $false, $true, 1, $false, $textD, 2) | Out-Null

# Docente(s) paragraph (9): old TEXT_C1 -> TEXT_A ; old TEXT_C2 -> TEXT_F
$r = $d.Paragraphs(9).Range
$r.Find.Execute($textC1, $true, $false, $false, $false, $false, $true, 1, $false, $textA, 2) | Out-Null
$r = $d.Paragraphs(9).Range
$r.Find.Execute($textC2, $true, $false, $false, $false, $false, $true, 1, $false, $textF, 2) | Out-Null

# Programa resumido paragraph (11): old TEXT_D -> TEXT_H
$r = $d.Paragraphs(11).Range
$r.Find.Execute($textD, $true, $false, $false, $false, $false, $true, 1, $false, $textH, 2) | Out-Null

# Programa paragraph (14): old TEXT_F -> TEXT_I
$r = $d.Paragraphs(14).Range
$r.Find.Execute($textF, $true, $false, $false, $false, $false, $true, 1, $false, $textI, 2) | Out-Null

# Avaliacao paragraph (17): old TEXT_H -> TEXT_J ; old TEXT_I -> TEXT_K ; old TEXT_J -> TEXT_C1
$r = $d.Paragraphs(17).Range
$r.Find.Execute($textH, $true, $false, $false, $false, $false, $true, 1, $false, $textJ, 2) | Out-Null
$r = $d.Paragraphs(17).Range
$r.Find.Execute($textI, $true, $false, $false, $false, $false, $true, 1, $false, $textK, 2) | Out-Null
$r = $d.Paragraphs(17).Range
$r.Find.Execute($textJ, $true, $false, $false, $false, $false, $true, 1, $false, $textC1, 2) | Out-Null

# Bibliografia paragraph (19): old TEXT_K -> TEXT_C2
$r = $d.Paragraphs(19).Range
$r.Find.Execute($textK, $true, $false, $false, $false, $false, $true, 1, $false, $textC2, 2) | Out-Null

# ---------------------------------------------------------------------
# Phase 3 - remove the two paragraphs that disappear entirely: the
# empty italic paragraph right after "Objetivos" (now paragraph 7,
# since paragraph 6 kept its place) and the italic English paragraph
# right after "Programa resumido" (now paragraph 12).
# ---------------------------------------------------------------------

$d.Paragraphs(7).Range.Delete() | Out-Null
$d.Paragraphs(11).Range.Delete() | Out-Null
